# Fix \n into doc_id. run again validation notebooks
#
# The "filename_pred" column (H) on Sheet1 pointed at the wrong
# predicted-source-document for a number of rows (an artifact of a
# newline that had leaked into a doc_id during parsing, which threw off
# the row/file alignment for several blocks of rows). This re-assigns the
# correct predicted filename for the affected rows. Two of those rows
# (17/18) also have their id_pred / id_match / id_match_label recomputed
# from the corrected filename_pred, so those are fixed up too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- filename_pred (column H) corrections ---------------------------------
$ws.Range("H6").Value2  = "documents-2024-11-30-1.json"
$ws.Range("H7").Value2  = "documents-2024-11-30-2.json"
$ws.Range("H8").Value2  = "documents-2024-12-01-1.json"
$ws.Range("H9").Value2  = "documents-2024-12-02-1.json"
$ws.Range("H10").Value2 = "documents-2024-12-02-5.json"
$ws.Range("H11").Value2 = "documents-2024-12-02-9.json"
$ws.Range("H12").Value2 = "documents-2024-12-02-8.json"

$ws.Range("H15").Value2 = "documents-2024-12-02-9.json"
$ws.Range("H16").Value2 = "documents-2024-12-02-8.json"

$ws.Range("H17").Value2 = "documents-2024-12-02-5.json"
$ws.Range("H18").Value2 = "documents-2024-12-02-4.json"

$ws.Range("H20").Value2 = "documents-2024-12-02-1.json"

$ws.Range("H24").Value2 = "documents-2024-12-02-5.json"
$ws.Range("H25").Value2 = "documents-2024-11-29-4.json"
$ws.Range("H26").Value2 = "documents-2024-11-30-2.json"
$ws.Range("H27").Value2 = "documents-2024-12-02-4.json"

$ws.Range("H30").Value2 = "documents-2024-12-02-1.json"
$ws.Range("H31").Value2 = "documents-2024-12-01-1.json"

# --- rows 17/18: id_pred swaps, so id_match / id_match_label flip ---------
$ws.Range("B17").Value2 = 4
$ws.Range("Q17").Value2 = $true
$ws.Range("R17").Value2 = "match"

$ws.Range("B18").Value2 = 3
$ws.Range("Q18").Value2 = $false
$ws.Range("R18").Value2 = "mismatch"
